$wb = $excel.ActiveWorkbook

# Helper-less approach: Excel COM auto-converts numeric-looking strings assigned
# to Range.Value into real numbers (losing exact text representation and the
# shared-string storage used by the workbook). Prefixing the value with a
# leading apostrophe forces Excel to store it as literal text (same as a user
# typing '0.123 into a cell); we then reset the cell style back to "Normal" so
# the quote-prefix formatting flag doesn't linger on the cell.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# NOTE: worksheet lookup by name is case-insensitive in this host, and this
# workbook has two sheets whose names differ only by case ("Vector_bf" vs.
# "Vector_BF"), so sheets are addressed by their (1-based) tab position to
# avoid ambiguity:
#   1 Funciones_Objetivo        2 Restricciones_del_lider
#   3 Restricciones_del_follower 4 Punto_modificado
#   5 Vector_bf                 6 Vector_BF
#   7 Vector_Alpha

# --- Restricciones_del_follower (sheet3) ---------------------------------
$ws3 = $wb.Worksheets.Item(3)

Set-TextValue $ws3.Range("B2") "-4.657691821664619"
Set-TextValue $ws3.Range("D2") "0.24011722556595838"
Set-TextValue $ws3.Range("E2") "0"
Set-TextValue $ws3.Range("F2") "0.054839693406650514"

Set-TextValue $ws3.Range("B3") "0.6576918216646188"
Set-TextValue $ws3.Range("D3") "0.1083236165390392"
Set-TextValue $ws3.Range("E3") "0.14589785305209468"

Set-TextValue $ws3.Range("B4") "-4.4622558915346"
Set-TextValue $ws3.Range("D4") "0.1102758390135593"
Set-TextValue $ws3.Range("E4") "0.4290202878062076"
Set-TextValue $ws3.Range("F4") "0.9920904723087388"

Set-TextValue $ws3.Range("B5") "3.6605518210954013"
Set-TextValue $ws3.Range("D5") "0.3168885247170169"
Set-TextValue $ws3.Range("E5") "0.8178656406623357"
Set-TextValue $ws3.Range("F5") "0"

Set-TextValue $ws3.Range("B6") "-9.066098643186933"
Set-TextValue $ws3.Range("D6") "0.4167665579899481"
Set-TextValue $ws3.Range("E6") "0"
Set-TextValue $ws3.Range("F6") "0.30309304447134744"

# --- Punto_modificado (sheet4) -------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

Set-TextValue $ws4.Range("A2") "5.875357499928848"
Set-TextValue $ws4.Range("B2") "4.657691821664619"
Set-TextValue $ws4.Range("C2") "2.112315956957238"

# --- Vector_bf (sheet5) ---------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-TextValue $ws5.Range("A2") "1.0970590052220022"

# --- Vector_BF (sheet6) ----------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-TextValue $ws6.Range("A2") "-4.68488454968627"
Set-TextValue $ws6.Range("A3") "-0.10886493274546893"
